# Update cryptos list with latest scraped prices and volume percentages
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.960.78"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").Value = "'1.878.59"
$ws.Range("E3").Value = "  -0.52%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("E5").Value = "  -3.72%  "
$ws.Range("D6").Value = "'242.67"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").Value = "'1.002"
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("D8").Value = "'0.3153"
$ws.Range("E8").Value = "  +1.13%  "
$ws.Range("D9").Value = "'0.07221"
$ws.Range("E9").Value = "  +0.67%  "
$ws.Range("D10").Value = "'24.70"
$ws.Range("E10").Value = "  -3.49%  "
$ws.Range("D11").Value = "'0.08376"
$ws.Range("E11").Value = "  -2.72%  "
$ws.Range("D12").Value = "'1.916.51"
$ws.Range("E12").Value = "  -5.50%  "
$ws.Range("D13").Value = "'0.7529"
$ws.Range("E13").Value = "  -1.28%  "
$ws.Range("D14").Value = "'5.433"
$ws.Range("E14").Value = "  +1.18%  "
$ws.Range("D15").Value = "'92.62"
$ws.Range("E15").Value = "  -1.00%  "
$ws.Range("D16").Value = "'29.989.38"
$ws.Range("E16").Value = "  +0.11%  "
$ws.Range("D17").Value = "'6.094"
$ws.Range("E17").Value = "  -1.33%  "
$ws.Range("D18").Value = "'251.60"
$ws.Range("E18").Value = "  +2.92%  "
$ws.Range("D20").Value = "'0.000007861"
$ws.Range("E20").Value = "  +0.54%  "
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("D22").Value = "'2.139.40"
$ws.Range("E22").Value = "  -2.77%  "
$ws.Range("D23").Value = "'8.057"
$ws.Range("E23").Value = "  +1.04%  "
$ws.Range("E24").Value = "  +0.21%  "
$ws.Range("D25").Value = "'0.1550"
$ws.Range("E25").Value = "  -5.99%  "
$ws.Range("D26").Value = "'9.268"
$ws.Range("E26").Value = "  -1.01%  "
$ws.Range("D27").Value = "'165.16"
$ws.Range("E27").Value = "  +1.62%  "
$ws.Range("D28").Value = "'18.72"
$ws.Range("E28").Value = "  -0.14%  "
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("D30").Value = "'1.514"
$ws.Range("E30").Value = "  +4.69%  "
$ws.Range("D31").Value = "'4.615"
$ws.Range("E31").Value = "  +2.23%  "
$ws.Range("E32").Value = "  +0.21%  "
$ws.Range("D33").Value = "'4.314"
$ws.Range("E33").Value = "  +5.20%  "
$ws.Range("D34").Value = "'0.05326"
$ws.Range("E34").Value = "  -1.89%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").Value = "'0.7482"
$ws.Range("E36").Value = "  +0.48%  "
$ws.Range("D37").Value = "'1.006"
$ws.Range("E37").Value = "  +0.29%  "
$ws.Range("D38").Value = "'2.693"
$ws.Range("E38").Value = "  -0.22%  "
$ws.Range("D39").Value = "'0.01968"
$ws.Range("E39").Value = "  +0.51%  "
$ws.Range("D40").Value = "'2.762"
$ws.Range("E40").Value = "  -0.74%  "
$ws.Range("D41").Value = "'0.4553"
$ws.Range("E41").Value = "  +2.01%  "
$ws.Range("D42").Value = "'1.113.00"
$ws.Range("E42").Value = "  +0.33%  "
$ws.Range("D43").Value = "'6.102"
$ws.Range("E43").Value = "  +0.34%  "
$ws.Range("D44").Value = "'72.49"
$ws.Range("E44").Value = "  -0.51%  "
$ws.Range("D45").Value = "'0.8538"
$ws.Range("E45").Value = "  +0.30%  "
$ws.Range("D46").Value = "'1.003"
$ws.Range("E46").Value = "  +0.20%  "
$ws.Range("D47").Value = "'103.85"
$ws.Range("E47").Value = "  +1.18%  "
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").Value = "'7.621"
$ws.Range("E48").Value = "  -0.07%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "'1.854"
$ws.Range("E49").Value = "  -0.63%  "
$ws.Range("D50").Value = "'2.038.66"
$ws.Range("E50").Value = "  -2.45%  "
$ws.Range("D51").Value = "'2.915"
$ws.Range("E51").Value = "  -2.10%  "
